# Apply value updates to team-specific transition matrix (Sheet1)
# per commit: "changes to team matrices from games pulled march 7"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1836065573770492
$ws.Range("C2").Value = 0.6
$ws.Range("J2").Value = 0.009836065573770493
$ws.Range("P2").Value = 0.1344262295081967
$ws.Range("S2").Value = 0.07213114754098361
$ws.Range("B3").Value = 0.005376344086021506
$ws.Range("C3").Value = 0.03225806451612903
$ws.Range("J3").Value = 0.02150537634408602
$ws.Range("P3").Value = 0.7741935483870968
$ws.Range("S3").Value = 0.1666666666666667
$ws.Range("B6").Value = 0.07589285714285714
$ws.Range("D6").Value = 0.008928571428571428
$ws.Range("F6").Value = 0.04910714285714286
$ws.Range("J6").Value = 0.3214285714285715
$ws.Range("O6").Value = 0.01339285714285714
$ws.Range("Q6").Value = 0.1607142857142857
$ws.Range("R6").Value = 0.07589285714285714
$ws.Range("S6").Value = 0.2946428571428572
$ws.Range("B7").Value = 0.07291666666666667
$ws.Range("D7").Value = 0.02083333333333333
$ws.Range("F7").Value = 0.05208333333333334
$ws.Range("J7").Value = 0.1302083333333333
$ws.Range("O7").Value = 0.02083333333333333
$ws.Range("Q7").Value = 0.2083333333333333
$ws.Range("R7").Value = 0.07291666666666667
$ws.Range("S7").Value = 0.421875
$ws.Range("B8").Value = 0.08043478260869565
$ws.Range("D8").Value = 0.02391304347826087
$ws.Range("E8").Value = 0.004347826086956522
$ws.Range("F8").Value = 0.04130434782608695
$ws.Range("J8").Value = 0.1065217391304348
$ws.Range("O8").Value = 0.0108695652173913
$ws.Range("Q8").Value = 0.2217391304347826
$ws.Range("R8").Value = 0.09782608695652174
$ws.Range("S8").Value = 0.4130434782608696
$ws.Range("B9").Value = 0.103448275862069
$ws.Range("D9").Value = 0.01477832512315271
$ws.Range("F9").Value = 0.05911330049261083
$ws.Range("J9").Value = 0.09359605911330049
$ws.Range("O9").Value = 0.009852216748768473
$ws.Range("Q9").Value = 0.2118226600985222
$ws.Range("R9").Value = 0.09852216748768473
$ws.Range("S9").Value = 0.4088669950738916
$ws.Range("B10").Value = 0.1166290443942814
$ws.Range("D10").Value = 0.0308502633559067
$ws.Range("E10").Value = 0.001504890895410083
$ws.Range("F10").Value = 0.06170052671181339
$ws.Range("J10").Value = 0.08502633559066967
$ws.Range("O10").Value = 0.01354401805869074
$ws.Range("Q10").Value = 0.2302483069977427
$ws.Range("R10").Value = 0.1075996990218209
$ws.Range("S10").Value = 0.3528969149736644
$ws.Range("G11").Value = 0.1258503401360544
$ws.Range("J11").Value = 0.1122448979591837
$ws.Range("K11").Value = 0.2006802721088435
$ws.Range("L11").Value = 0.5408163265306123
$ws.Range("S11").Value = 0.02040816326530612
$ws.Range("G12").Value = 0.7283950617283951
$ws.Range("J12").Value = 0.228395061728395
$ws.Range("L12").Value = 0.0308641975308642
$ws.Range("S12").Value = 0.01234567901234568
$ws.Range("G13").Value = 0.7843137254901961
$ws.Range("J13").Value = 0.196078431372549
$ws.Range("S13").Value = 0.0196078431372549
$ws.Range("F15").Value = 0.02643171806167401
$ws.Range("H15").Value = 0.118942731277533
$ws.Range("I15").Value = 0.1013215859030837
$ws.Range("J15").Value = 0.4537444933920705
$ws.Range("K15").Value = 0.06167400881057269
$ws.Range("M15").Value = 0.013215859030837
$ws.Range("N15").Value = 0.004405286343612335
$ws.Range("O15").Value = 0.06167400881057269
$ws.Range("S15").Value = 0.1585903083700441
$ws.Range("F16").Value = 0.004444444444444444
$ws.Range("H16").Value = 0.2355555555555555
$ws.Range("I16").Value = 0.05777777777777778
$ws.Range("J16").Value = 0.44
$ws.Range("K16").Value = 0.08888888888888889
$ws.Range("M16").Value = 0.03111111111111111
$ws.Range("O16").Value = 0.04444444444444445
$ws.Range("S16").Value = 0.09777777777777778
$ws.Range("F17").Value = 0.03250478011472276
$ws.Range("H17").Value = 0.2007648183556405
$ws.Range("I17").Value = 0.08413001912045889
$ws.Range("J17").Value = 0.4130019120458891
$ws.Range("K17").Value = 0.08795411089866156
$ws.Range("M17").Value = 0.02103250478011472
$ws.Range("N17").Value = 0.001912045889101338
$ws.Range("O17").Value = 0.07839388145315487
$ws.Range("S17").Value = 0.08030592734225621
$ws.Range("F18").Value = 0.01694915254237288
$ws.Range("H18").Value = 0.1694915254237288
$ws.Range("I18").Value = 0.08050847457627118
$ws.Range("J18").Value = 0.4830508474576271
$ws.Range("K18").Value = 0.06779661016949153
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.0635593220338983
$ws.Range("S18").Value = 0.1016949152542373
$ws.Range("F19").Value = 0.02892561983471074
$ws.Range("H19").Value = 0.1950413223140496
$ws.Range("I19").Value = 0.0884297520661157
$ws.Range("J19").Value = 0.3743801652892562
$ws.Range("K19").Value = 0.112396694214876
$ws.Range("M19").Value = 0.02231404958677686
$ws.Range("O19").Value = 0.07107438016528926
$ws.Range("S19").Value = 0.1074380165289256
